# Rename the "AddressBook" sample app to "ThanePark" throughout the
# UndoRedo sequence diagram on slide 1 (docs/diagrams/UndoRedoSequenceDiagram.pptx).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Lifeline header box ":AddressBookParser" (wrapped over two lines as
#    ":Address" / "BookParser") -> ":ThaneParkParser" (wrapped as
#    ":ThanePark" / "Parser").
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange
if ($tr.Text -ne (":Address" + [char]13 + "BookParser")) {
    throw "Shape 6: unexpected text [$($tr.Text)]"
}
# Drop "Book" (the text right after the paragraph break) first ...
$book = $tr.Characters(10, 4)
if ($book.Text -ne "Book") { throw "Shape 6: expected 'Book', got [$($book.Text)]" }
$book.Text = ""
# ... then turn "Address" into "ThanePark".
$addr = $tr.Characters(2, 7)
if ($addr.Text -ne "Address") { throw "Shape 6: expected 'Address', got [$($addr.Text)]" }
$addr.Text = "ThanePark"
if ($tr.Text -ne (":ThanePark" + [char]13 + "Parser")) {
    throw "Shape 6: unexpected result [$($tr.Text)]"
}

# 2) "undoAddressBook()" call label -> "undoThanePark()"
$sh = $s.Shapes.Item(19)
$tr = $sh.TextFrame.TextRange
if ($tr.Text -ne "undoAddressBook()") { throw "Shape 19: unexpected text [$($tr.Text)]" }
$run = $tr.Characters(5, 11)
if ($run.Text -ne "AddressBook") { throw "Shape 19: expected 'AddressBook', got [$($run.Text)]" }
$run.Text = "ThanePark"
if ($tr.Text -ne "undoThanePark()") { throw "Shape 19: unexpected result [$($tr.Text)]" }

# 3) ":VersionedAddressBook" lifeline -> ":VersionedThanePark"
$sh = $s.Shapes.Item(23)
$tr = $sh.TextFrame.TextRange
if ($tr.Text -ne ":VersionedAddressBook") { throw "Shape 23: unexpected text [$($tr.Text)]" }
$run = $tr.Characters(2, 20)
if ($run.Text -ne "VersionedAddressBook") { throw "Shape 23: expected 'VersionedAddressBook', got [$($run.Text)]" }
$run.Text = "VersionedThanePark"
if ($tr.Text -ne ":VersionedThanePark") { throw "Shape 23: unexpected result [$($tr.Text)]" }

# 4) "resetData(ReadOnlyAddressBook)" call label -> "resetData(ReadOnlyThanePark)"
$sh = $s.Shapes.Item(35)
$tr = $sh.TextFrame.TextRange
if ($tr.Text -ne "resetData(ReadOnlyAddressBook)") { throw "Shape 35: unexpected text [$($tr.Text)]" }
$run = $tr.Characters(11, 19)
if ($run.Text -ne "ReadOnlyAddressBook") { throw "Shape 35: expected 'ReadOnlyAddressBook', got [$($run.Text)]" }
$run.Text = "ReadOnlyThanePark"
if ($tr.Text -ne "resetData(ReadOnlyThanePark)") { throw "Shape 35: unexpected result [$($tr.Text)]" }

Write-Output "Renamed AddressBook -> ThanePark in 4 shapes on slide 1"
